$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 605.0952
$ws.Range("I4").Value = 515.35297
$ws.Range("K4").Value = 515.35297
$ws.Range("M4").Value = -401.35297
$ws.Range("H5").Value = 210.38461
$ws.Range("I5").Value = 140.36363
$ws.Range("K5").Value = 140.36363
$ws.Range("M5").Value = -25.36363
$ws.Range("H17").Value = 1194.6
$ws.Range("J17").Value = 1194.6
$ws.Range("L17").Value = 3583.8
$ws.Range("N17").Value = -3919.8
$ws.Range("H40").Value = 71431016
$ws.Range("I40").Value = 2309.75
$ws.Range("K40").Value = 2309.75
$ws.Range("M40").Value = -2134.75
$ws.Range("H100").Value = 4098.815
$ws.Range("I100").Value = 1256.4286
$ws.Range("J100").Value = 7159.846
$ws.Range("K100").Value = 1256.4286
$ws.Range("L100").Value = 7159.846
$ws.Range("M100").Value = -715.4286
$ws.Range("N100").Value = -8241.846
$ws.Range("H112").Value = 1789869.1
$ws.Range("I112").Value = 712
$ws.Range("J112").Value = 2783845.2
$ws.Range("K112").Value = 2136
$ws.Range("L112").Value = 8351535.600000001
$ws.Range("M112").Value = -1028
$ws.Range("N112").Value = -8353751.600000001
$ws.Range("H132").Value = 4822.2812
$ws.Range("I132").Value = 3236.6667
$ws.Range("K132").Value = 9710.000100000001
$ws.Range("M132").Value = -7180.000100000001
$ws.Range("H135").Value = 1432.0938
$ws.Range("I135").Value = 1278.138
$ws.Range("J135").Value = 2920.3333
$ws.Range("K135").Value = 11503.242
$ws.Range("L135").Value = 26282.9997
$ws.Range("M135").Value = -8968.241999999998
$ws.Range("N135").Value = -31352.9997
$ws.Range("H138").Value = 5139.4473
$ws.Range("J138").Value = 7678.125
$ws.Range("L138").Value = 23034.375
$ws.Range("N138").Value = -33314.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4841.143
$ws.Range("I45").Value = 4347
$ws.Range("J45").Value = 5500
$ws.Range("K45").Value = 4347
$ws.Range("L45").Value = 5500
$ws.Range("M45").Value = -3970
$ws.Range("N45").Value = -6254
$ws.Range("H61").Value = 980657.7
$ws.Range("I61").Value = 35081.25
$ws.Range("K61").Value = 35081.25
$ws.Range("M61").Value = -34869.25
$ws.Range("H122").Value = 6616.3335
$ws.Range("I122").Value = 5226.615
$ws.Range("K122").Value = 15679.845
$ws.Range("M122").Value = -13229.845
$ws.Range("H132").Value = 6669734
$ws.Range("I132").Value = 3456.6365
$ws.Range("K132").Value = 10369.9095
$ws.Range("M132").Value = -7839.9095
$ws.Range("H136").Value = 980657.7
$ws.Range("I136").Value = 35081.25
$ws.Range("K136").Value = 105243.75
$ws.Range("M136").Value = -102693.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 820603.1
$ws.Range("I105").Value = 1635672.5
$ws.Range("J105").Value = 5533.7856
$ws.Range("K105").Value = 1635672.5
$ws.Range("L105").Value = 5533.7856
$ws.Range("M105").Value = -1633925.5
$ws.Range("N105").Value = -9027.785599999999
$ws.Range("H134").Value = 50002500
$ws.Range("I134").Value = 5000
$ws.Range("K134").Value = 15000
$ws.Range("M134").Value = -12465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1363.6
$ws.Range("I16").Value = 1374.25
$ws.Range("J16").Value = 1321
$ws.Range("K16").Value = 1374.25
$ws.Range("L16").Value = 1321
$ws.Range("M16").Value = -1087.25
$ws.Range("N16").Value = -1895
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H31").Value = 2454.7188
$ws.Range("I31").Value = 2251.7273
$ws.Range("J31").Value = 2901.3
$ws.Range("K31").Value = 2251.7273
$ws.Range("L31").Value = 2901.3
$ws.Range("M31").Value = -1956.7273
$ws.Range("N31").Value = -3491.3
$ws.Range("H34").Value = 2454.7188
$ws.Range("I34").Value = 2251.7273
$ws.Range("J34").Value = 2901.3
$ws.Range("K34").Value = 2251.7273
$ws.Range("L34").Value = 2901.3
$ws.Range("M34").Value = -2049.7273
$ws.Range("N34").Value = -3305.3
$ws.Range("H47").Value = 207010.5
$ws.Range("I47").Value = 357688
$ws.Range("J47").Value = 56333
$ws.Range("K47").Value = 357688
$ws.Range("L47").Value = 56333
$ws.Range("M47").Value = -357122
$ws.Range("N47").Value = -57465
$ws.Range("H98").Value = 99998.5
$ws.Range("J98").Value = 99998.5
$ws.Range("L98").Value = 99998.5
$ws.Range("N98").Value = -104490.5
$ws.Range("H100").Value = 99998.664
$ws.Range("J100").Value = 99998.664
$ws.Range("L100").Value = 99998.664
$ws.Range("N100").Value = -102162.664
$ws.Range("H113").Value = 1363.6
$ws.Range("I113").Value = 1374.25
$ws.Range("J113").Value = 1321
$ws.Range("K113").Value = 1374.25
$ws.Range("L113").Value = 1321
$ws.Range("M113").Value = 795.75
$ws.Range("N113").Value = -5661
$ws.Range("H132").Value = 2586.6667
$ws.Range("I132").Value = 2586.6667
$ws.Range("K132").Value = 7760.000100000001
$ws.Range("M132").Value = -5230.000100000001
$ws.Range("H134").Value = 5422.5
$ws.Range("I134").Value = 5422.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 16267.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -13732.5
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1410.825
$ws.Range("J22").Value = 1509.25
$ws.Range("L22").Value = 4527.75
$ws.Range("N22").Value = -4865.75
$ws.Range("H27").Value = 1410.825
$ws.Range("J27").Value = 1509.25
$ws.Range("L27").Value = 4527.75
$ws.Range("N27").Value = -4731.75
$ws.Range("H68").Value = 1763.909
$ws.Range("J68").Value = 1862.875
$ws.Range("L68").Value = 5588.625
$ws.Range("N68").Value = -7210.625
$ws.Range("H71").Value = 1763.909
$ws.Range("J71").Value = 1862.875
$ws.Range("L71").Value = 16765.875
$ws.Range("N71").Value = -24877.875
$ws.Range("H74").Value = 14186
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 14186
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 42558
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -44680
$ws.Range("H77").Value = 14186
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 14186
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 127674
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -138282
$ws.Range("H98").Value = 690.8
$ws.Range("J98").Value = 683.75
$ws.Range("L98").Value = 2051.25
$ws.Range("N98").Value = -5047.25
$ws.Range("H139").Value = 2657.3704
$ws.Range("I139").Value = 1654.25
$ws.Range("J139").Value = 3079.7368
$ws.Range("K139").Value = 4962.75
$ws.Range("L139").Value = 9239.2104
$ws.Range("M139").Value = 177.25
$ws.Range("N139").Value = -19519.2104

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9005.4
$ws.Range("I70").Value = 8595.333000000001
$ws.Range("J70").Value = 9181.143
$ws.Range("K70").Value = 8595.333000000001
$ws.Range("L70").Value = 9181.143
$ws.Range("M70").Value = -8325.333000000001
$ws.Range("N70").Value = -9721.143
$ws.Range("H73").Value = 9005.4
$ws.Range("I73").Value = 8595.333000000001
$ws.Range("J73").Value = 9181.143
$ws.Range("K73").Value = 8595.333000000001
$ws.Range("L73").Value = 9181.143
$ws.Range("M73").Value = -7659.333000000001
$ws.Range("N73").Value = -11053.143
$ws.Range("H80").Value = 2483.5715
$ws.Range("I80").Value = 1531
$ws.Range("K80").Value = 1531
$ws.Range("M80").Value = -533
$ws.Range("H83").Value = 2483.5715
$ws.Range("I83").Value = 1531
$ws.Range("K83").Value = 7655
$ws.Range("M83").Value = -2663
$ws.Range("H122").Value = 1959
$ws.Range("I122").Value = 1959
$ws.Range("K122").Value = 5877
$ws.Range("M122").Value = -3427
$ws.Range("H126").Value = 2409.6
$ws.Range("I126").Value = 2115
$ws.Range("J126").Value = 4324.5
$ws.Range("K126").Value = 6345
$ws.Range("L126").Value = 12973.5
$ws.Range("M126").Value = -3875
$ws.Range("N126").Value = -17913.5
$ws.Range("H132").Value = 9549276
$ws.Range("I132").Value = 4612.8184
$ws.Range("K132").Value = 13838.4552
$ws.Range("M132").Value = -11308.4552
$ws.Range("H141").Value = 79992.664
$ws.Range("J141").Value = 79992.664
$ws.Range("L141").Value = 79992.664
$ws.Range("N141").Value = -90352.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1452.6875
$ws.Range("I55").Value = 1391.3334
$ws.Range("J55").Value = 1489.5
$ws.Range("K55").Value = 1391.3334
$ws.Range("L55").Value = 1489.5
$ws.Range("M55").Value = -1218.3334
$ws.Range("N55").Value = -1835.5
$ws.Range("H122").Value = 3631.76
$ws.Range("I122").Value = 2886.3333
$ws.Range("J122").Value = 4749.9
$ws.Range("K122").Value = 8658.999899999999
$ws.Range("L122").Value = 14249.7
$ws.Range("M122").Value = -6208.999899999999
$ws.Range("N122").Value = -19149.7
$ws.Range("H132").Value = 2189.1428
$ws.Range("I132").Value = 1982.8286
$ws.Range("J132").Value = 2704.9285
$ws.Range("K132").Value = 5948.4858
$ws.Range("L132").Value = 8114.7855
$ws.Range("M132").Value = -3418.4858
$ws.Range("N132").Value = -13174.7855
$ws.Range("H136").Value = 2244.4524
$ws.Range("I136").Value = 2142.9688
$ws.Range("K136").Value = 6428.9064
$ws.Range("M136").Value = -3878.9064

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 836174.0600000001
$ws.Range("I132").Value = 3208.9
$ws.Range("J132").Value = 5001000
$ws.Range("K132").Value = 9626.700000000001
$ws.Range("L132").Value = 15003000
$ws.Range("M132").Value = -7096.700000000001
$ws.Range("N132").Value = -15008060
$ws.Range("H141").Value = 115323.914
$ws.Range("J141").Value = 126523.7
$ws.Range("L141").Value = 126523.7
$ws.Range("N141").Value = -136883.7
